$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: fill in the Data Structure / Difficulty / From columns for "Number of Islands"
$ws.Range("D36").Value = "Graph/UF"
$ws.Range("D36").HorizontalAlignment = -4108
$ws.Range("E36").Value = "medium"
$ws.Range("E36").HorizontalAlignment = -4108
$ws.Range("F36").Value = "leetcode 200"
$ws.Range("F36").HorizontalAlignment = -4108

# Row 37: new entry "Number of Islands II"
$ws.Range("A37").Value = 36
$ws.Range("A37").HorizontalAlignment = -4108
$ws.Range("B37").Value = "Number of Islands II"
$ws.Range("B37").HorizontalAlignment = -4131
$ws.Range("D37").Value = "Graph/UF"
$ws.Range("D37").HorizontalAlignment = -4108
$ws.Range("E37").Value = "medium"
$ws.Range("E37").HorizontalAlignment = -4108
$ws.Range("F37").Value = "leetcode 305"
$ws.Range("F37").HorizontalAlignment = -4108

# Update the view: selection moves to F37, no frozen/top-left offset
$ws.Range("F37").Select()
